# Apply KHL probabilities tour update
# - Updates existing match rows 2-4 with refreshed odds/probabilities
# - Adds two new match rows (5 and 6) for additional games

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2,2).Value = '2025-12-21T13:30:00'
$ws.Cells.Item(2,4).Value = 'Амур'
$ws.Cells.Item(2,5).Value = 897889
$ws.Cells.Item(2,6).Value = 'https://text.khl.ru/text/897889.html'
$ws.Cells.Item(2,7).Value = 2.119583
$ws.Cells.Item(2,8).Value = 1.566641
$ws.Cells.Item(2,9).Value = 3.547501
$ws.Cells.Item(2,10).Value = 4.529412
$ws.Cells.Item(2,11).Value = 3.324498
$ws.Cells.Item(2,12).Value = 2.557071
$ws.Cells.Item(2,13).Value = 3.686225
$ws.Cells.Item(2,14).Value = 30.079763
$ws.Cells.Item(2,15).Value = 27.002774
$ws.Cells.Item(2,16).Value = 57.082537
$ws.Cells.Item(2,17).Value = -0.012706
$ws.Cells.Item(2,18).Value = -0.133671
$ws.Cells.Item(2,19).Value = 0.540635
$ws.Cells.Item(2,20).Value = 0.160882
$ws.Cells.Item(2,21).Value = 0.297732
$ws.Cells.Item(2,22).Value = 0.162089
$ws.Cells.Item(2,23).Value = 0.837159
$ws.Cells.Item(2,24).Value = 0.301222
$ws.Cells.Item(2,25).Value = 0.698026
$ws.Cells.Item(2,26).Value = 0.464885
$ws.Cells.Item(2,27).Value = 0.534363
$ws.Cells.Item(2,28).Value = 0.625318
$ws.Cells.Item(2,29).Value = 0.37393
$ws.Cells.Item(2,30).Value = 0.760118
$ws.Cells.Item(2,31).Value = 0.23913
$ws.Cells.Item(2,32).Value = 0.844359
$ws.Cells.Item(2,33).Value = 0.155641
$ws.Cells.Item(2,34).Value = 0.64547
$ws.Cells.Item(2,35).Value = 0.35453
$ws.Cells.Item(2,36).Value = 0.7242150000000001
$ws.Cells.Item(2,37).Value = 0.275785
$ws.Cells.Item(2,38).Value = 0.470741
$ws.Cells.Item(2,39).Value = 0.529259
$ws.Cells.Item(2,40).Value = 0.829865
$ws.Cells.Item(2,41).Value = 0.625482

# Row 3
$ws.Cells.Item(3,2).Value = '2025-12-21T15:00:00'
$ws.Cells.Item(3,3).Value = 'Автомобилист'
$ws.Cells.Item(3,4).Value = 'Адмирал'
$ws.Cells.Item(3,5).Value = 897887
$ws.Cells.Item(3,6).Value = 'https://text.khl.ru/text/897887.html'
$ws.Cells.Item(3,7).Value = 4.694118
$ws.Cells.Item(3,8).Value = 1.233333
$ws.Cells.Item(3,9).Value = 3.458824
$ws.Cells.Item(3,10).Value = 3.167914
$ws.Cells.Item(3,11).Value = 3.931016
$ws.Cells.Item(3,12).Value = 2.346078
$ws.Cells.Item(3,13).Value = 5.927451
$ws.Cells.Item(3,14).Value = 31.854488
$ws.Cells.Item(3,15).Value = 26.681133
$ws.Cells.Item(3,16).Value = 58.535621
$ws.Cells.Item(3,18).Value = -0.16
$ws.Cells.Item(3,19).Value = 0.663741
$ws.Cells.Item(3,20).Value = 0.135116
$ws.Cells.Item(3,21).Value = 0.198614
$ws.Cells.Item(3,22).Value = 0.128137
$ws.Cells.Item(3,23).Value = 0.869335
$ws.Cells.Item(3,24).Value = 0.249676
$ws.Cells.Item(3,25).Value = 0.747796
$ws.Cells.Item(3,26).Value = 0.402259
$ws.Cells.Item(3,27).Value = 0.595213
$ws.Cells.Item(3,28).Value = 0.5618880000000001
$ws.Cells.Item(3,29).Value = 0.435584
$ws.Cells.Item(3,30).Value = 0.705031
$ws.Cells.Item(3,31).Value = 0.29244
$ws.Cells.Item(3,32).Value = 0.903235
$ws.Cells.Item(3,33).Value = 0.096765
$ws.Cells.Item(3,34).Value = 0.751613
$ws.Cells.Item(3,35).Value = 0.248387
$ws.Cells.Item(3,36).Value = 0.679633
$ws.Cells.Item(3,37).Value = 0.320367
$ws.Cells.Item(3,38).Value = 0.416142
$ws.Cells.Item(3,39).Value = 0.583858
$ws.Cells.Item(3,40).Value = 0.894209
$ws.Cells.Item(3,41).Value = 0.493499

# Row 4
$ws.Cells.Item(4,2).Value = '2025-12-21T15:00:00'
$ws.Cells.Item(4,3).Value = 'Барыс'
$ws.Cells.Item(4,4).Value = 'ХК Сочи'
$ws.Cells.Item(4,5).Value = 897888
$ws.Cells.Item(4,6).Value = 'https://text.khl.ru/text/897888.html'
$ws.Cells.Item(4,7).Value = 1.657143
$ws.Cells.Item(4,8).Value = 1.2
$ws.Cells.Item(4,9).Value = 2.613665
$ws.Cells.Item(4,10).Value = 5.95
$ws.Cells.Item(4,11).Value = 3.803571
$ws.Cells.Item(4,12).Value = 1.906832
$ws.Cells.Item(4,13).Value = 2.857143
$ws.Cells.Item(4,14).Value = 26.696492
$ws.Cells.Item(4,15).Value = 24.465028
$ws.Cells.Item(4,16).Value = 51.161521
$ws.Cells.Item(4,17).Value = -0.16
$ws.Cells.Item(4,19).Value = 0.7183659999999999
$ws.Cells.Item(4,20).Value = 0.127573
$ws.Cells.Item(4,21).Value = 0.152113
$ws.Cells.Item(4,22).Value = 0.178976
$ws.Cells.Item(4,23).Value = 0.819076
$ws.Cells.Item(4,24).Value = 0.325686
$ws.Cells.Item(4,25).Value = 0.672366
$ws.Cells.Item(4,26).Value = 0.493241
$ws.Cells.Item(4,27).Value = 0.504811
$ws.Cells.Item(4,28).Value = 0.652708
$ws.Cells.Item(4,29).Value = 0.345344
$ws.Cells.Item(4,30).Value = 0.782797
$ws.Cells.Item(4,31).Value = 0.215255
$ws.Cells.Item(4,32).Value = 0.892923
$ws.Cells.Item(4,33).Value = 0.107077
$ws.Cells.Item(4,34).Value = 0.73168
$ws.Cells.Item(4,35).Value = 0.26832
$ws.Cells.Item(4,36).Value = 0.56819
$ws.Cells.Item(4,37).Value = 0.43181
$ws.Cells.Item(4,38).Value = 0.298125
$ws.Cells.Item(4,39).Value = 0.701875
$ws.Cells.Item(4,40).Value = 0.927381
$ws.Cells.Item(4,41).Value = 0.442138

# Row 5
$ws.Cells.Item(5,1).Value = 1369
$ws.Cells.Item(5,2).Value = '2025-12-21T17:00:00'
$ws.Cells.Item(5,3).Value = 'ЦСКА'
$ws.Cells.Item(5,4).Value = 'Динамо М'
$ws.Cells.Item(5,5).Value = 897891
$ws.Cells.Item(5,6).Value = 'https://text.khl.ru/text/897891.html'
$ws.Cells.Item(5,7).Value = 3.097291
$ws.Cells.Item(5,8).Value = 4.212037
$ws.Cells.Item(5,9).Value = 1.517647
$ws.Cells.Item(5,10).Value = 1.466667
$ws.Cells.Item(5,11).Value = 2.281979
$ws.Cells.Item(5,12).Value = 2.864842
$ws.Cells.Item(5,13).Value = 7.309328
$ws.Cells.Item(5,14).Value = 27.052298
$ws.Cells.Item(5,15).Value = 31.127614
$ws.Cells.Item(5,16).Value = 58.179911
$ws.Cells.Item(5,17).Value = 0.012972
$ws.Cells.Item(5,18).Value = 0.09715799999999999
$ws.Cells.Item(5,19).Value = 0.313238
$ws.Cells.Item(5,20).Value = 0.175449
$ws.Cells.Item(5,21).Value = 0.511087
$ws.Cells.Item(5,22).Value = 0.245018
$ws.Cells.Item(5,23).Value = 0.754756
$ws.Cells.Item(5,24).Value = 0.41512
$ws.Cells.Item(5,25).Value = 0.584654
$ws.Cells.Item(5,26).Value = 0.590217
$ws.Cells.Item(5,27).Value = 0.409557
$ws.Cells.Item(5,28).Value = 0.740415
$ws.Cells.Item(5,29).Value = 0.259359
$ws.Cells.Item(5,30).Value = 0.85085
$ws.Cells.Item(5,31).Value = 0.148924
$ws.Cells.Item(5,32).Value = 0.664969
$ws.Cells.Item(5,33).Value = 0.335031
$ws.Cells.Item(5,34).Value = 0.399177
$ws.Cells.Item(5,35).Value = 0.600823
$ws.Cells.Item(5,36).Value = 0.779734
$ws.Cells.Item(5,37).Value = 0.220266
$ws.Cells.Item(5,38).Value = 0.545858
$ws.Cells.Item(5,39).Value = 0.454142
$ws.Cells.Item(5,40).Value = 0.664817
$ws.Cells.Item(5,41).Value = 0.826832

# Row 6
$ws.Cells.Item(6,1).Value = 1369
$ws.Cells.Item(6,2).Value = '2025-12-21T17:10:00'
$ws.Cells.Item(6,3).Value = 'Динамо Мн'
$ws.Cells.Item(6,4).Value = 'Лада'
$ws.Cells.Item(6,5).Value = 897890
$ws.Cells.Item(6,6).Value = 'https://text.khl.ru/text/897890.html'
$ws.Cells.Item(6,7).Value = 5.6
$ws.Cells.Item(6,8).Value = 1.44
$ws.Cells.Item(6,9).Value = 1
$ws.Cells.Item(6,10).Value = 4.161709
$ws.Cells.Item(6,11).Value = 4.880855
$ws.Cells.Item(6,12).Value = 1.22
$ws.Cells.Item(6,13).Value = 7.04
$ws.Cells.Item(6,14).Value = 41.531495
$ws.Cells.Item(6,15).Value = 21.68186
$ws.Cells.Item(6,16).Value = 63.213355
$ws.Cells.Item(6,17).Value = 0.16
$ws.Cells.Item(6,18).Value = -0.16
$ws.Cells.Item(6,19).Value = 0.897023
$ws.Cells.Item(6,20).Value = 0.054868
$ws.Cells.Item(6,21).Value = 0.036448
$ws.Cells.Item(6,22).Value = 0.142428
$ws.Cells.Item(6,23).Value = 0.845911
$ws.Cells.Item(6,24).Value = 0.271784
$ws.Cells.Item(6,25).Value = 0.716556
$ws.Cells.Item(6,26).Value = 0.429619
$ws.Cells.Item(6,27).Value = 0.55872
$ws.Cells.Item(6,28).Value = 0.590108
$ws.Cells.Item(6,29).Value = 0.398232
$ws.Cells.Item(6,30).Value = 0.729982
$ws.Cells.Item(6,31).Value = 0.258358
$ws.Cells.Item(6,32).Value = 0.955361
$ws.Cells.Item(6,33).Value = 0.044639
$ws.Cells.Item(6,34).Value = 0.8649480000000001
$ws.Cells.Item(6,35).Value = 0.135052
$ws.Cells.Item(6,36).Value = 0.344589
$ws.Cells.Item(6,37).Value = 0.655411
$ws.Cells.Item(6,38).Value = 0.124879
$ws.Cells.Item(6,39).Value = 0.875121
$ws.Cells.Item(6,40).Value = 0.97632
$ws.Cells.Item(6,41).Value = 0.189048
